$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in rows 58 and 59, columns C:F (0 -> 1)
$ws.Range("C58:F58").Value = 1
$ws.Range("C59:F59").Value = 1

# Update the view: scroll to A40 and select F62
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F62").Select()
